$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 44400
$ws.Cells.Item(2, 8).Value = "Española"
$ws.Cells.Item(2, 9).Value = "Primera"
$ws.Cells.Item(2, 10).Value = 70
$ws.Cells.Item(2, 11).Value = 15000
$ws.Cells.Item(2, 12).Value = 15000
$ws.Cells.Item(2, 13).Value = 15000
$ws.Cells.Item(2, 14).Value = "`$/caja 30 unidades"
$ws.Cells.Item(2, 16).Value = 500
$ws.Cells.Item(2, 17).Value = 30

# Row 3
$ws.Cells.Item(3, 4).Value = 44495
$ws.Cells.Item(3, 8).Value = "Madrigal"
$ws.Cells.Item(3, 9).Value = "Primera"
$ws.Cells.Item(3, 10).Value = 130
$ws.Cells.Item(3, 11).Value = 11000
$ws.Cells.Item(3, 12).Value = 11000
$ws.Cells.Item(3, 13).Value = 11000
$ws.Cells.Item(3, 14).Value = "`$/caja 40 unidades"
$ws.Cells.Item(3, 16).Value = 275
$ws.Cells.Item(3, 17).Value = 40

# Row 4
$ws.Cells.Item(4, 4).Value = 44162
$ws.Cells.Item(4, 8).Value = "Madrigal"
$ws.Cells.Item(4, 9).Value = "Primera"
$ws.Cells.Item(4, 10).Value = 50
$ws.Cells.Item(4, 11).Value = 10000
$ws.Cells.Item(4, 12).Value = 10000
$ws.Cells.Item(4, 13).Value = 10000
$ws.Cells.Item(4, 14).Value = "`$/caja 40 unidades"
$ws.Cells.Item(4, 16).Value = 250
$ws.Cells.Item(4, 17).Value = 40

# Row 5
$ws.Cells.Item(5, 4).Value = 44446
$ws.Cells.Item(5, 8).Value = "Madrigal"
$ws.Cells.Item(5, 9).Value = "Primera"
$ws.Cells.Item(5, 10).Value = 120
$ws.Cells.Item(5, 11).Value = 16000
$ws.Cells.Item(5, 12).Value = 16000
$ws.Cells.Item(5, 13).Value = 16000
$ws.Cells.Item(5, 14).Value = "`$/caja 40 unidades"
$ws.Cells.Item(5, 16).Value = 400
$ws.Cells.Item(5, 17).Value = 40

# Row 6
$ws.Cells.Item(6, 4).Value = 44161
$ws.Cells.Item(6, 8).Value = "Madrigal"
$ws.Cells.Item(6, 9).Value = "Primera"
$ws.Cells.Item(6, 10).Value = 30
$ws.Cells.Item(6, 11).Value = 11000
$ws.Cells.Item(6, 12).Value = 11000
$ws.Cells.Item(6, 13).Value = 11000
$ws.Cells.Item(6, 14).Value = "`$/caja 40 unidades"
$ws.Cells.Item(6, 16).Value = 275
$ws.Cells.Item(6, 17).Value = 40

# Row 7
$ws.Cells.Item(7, 4).Value = 44166
$ws.Cells.Item(7, 8).Value = "Madrigal"
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 80
$ws.Cells.Item(7, 11).Value = 10000
$ws.Cells.Item(7, 12).Value = 10000
$ws.Cells.Item(7, 13).Value = 10000
$ws.Cells.Item(7, 14).Value = "`$/caja 40 unidades"
$ws.Cells.Item(7, 16).Value = 250
$ws.Cells.Item(7, 17).Value = 40

# Row 8
$ws.Cells.Item(8, 4).Value = 44491
$ws.Cells.Item(8, 8).Value = "Madrigal"
$ws.Cells.Item(8, 9).Value = "Primera"
$ws.Cells.Item(8, 10).Value = 200
$ws.Cells.Item(8, 11).Value = 11000
$ws.Cells.Item(8, 12).Value = 11000
$ws.Cells.Item(8, 13).Value = 11000
$ws.Cells.Item(8, 14).Value = "`$/caja 40 unidades"
$ws.Cells.Item(8, 16).Value = 275
$ws.Cells.Item(8, 17).Value = 40

# Row 9
$ws.Cells.Item(9, 4).Value = 44176
$ws.Cells.Item(9, 8).Value = "Madrigal"
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 80
$ws.Cells.Item(9, 11).Value = 11000
$ws.Cells.Item(9, 12).Value = 11000
$ws.Cells.Item(9, 13).Value = 11000
$ws.Cells.Item(9, 14).Value = "`$/caja 40 unidades"
$ws.Cells.Item(9, 16).Value = 275
$ws.Cells.Item(9, 17).Value = 40

# Row 10
$ws.Cells.Item(10, 4).Value = 44390
$ws.Cells.Item(10, 8).Value = "Española"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 80
$ws.Cells.Item(10, 11).Value = 16000
$ws.Cells.Item(10, 12).Value = 16000
$ws.Cells.Item(10, 13).Value = 16000
$ws.Cells.Item(10, 14).Value = "`$/caja 30 unidades"
$ws.Cells.Item(10, 16).Value = 533
$ws.Cells.Item(10, 17).Value = 30

# Row 11
$ws.Cells.Item(11, 4).Value = 44481
$ws.Cells.Item(11, 8).Value = "Madrigal"
$ws.Cells.Item(11, 9).Value = "Segunda"
$ws.Cells.Item(11, 10).Value = 120
$ws.Cells.Item(11, 11).Value = 11000
$ws.Cells.Item(11, 12).Value = 11000
$ws.Cells.Item(11, 13).Value = 11000
$ws.Cells.Item(11, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(11, 16).Value = 220
$ws.Cells.Item(11, 17).Value = 50

# Row 12
$ws.Cells.Item(12, 4).Value = 44407
$ws.Cells.Item(12, 8).Value = "Española"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 100
$ws.Cells.Item(12, 11).Value = 18000
$ws.Cells.Item(12, 12).Value = 18000
$ws.Cells.Item(12, 13).Value = 18000
$ws.Cells.Item(12, 14).Value = "`$/caja 30 unidades"
$ws.Cells.Item(12, 16).Value = 600
$ws.Cells.Item(12, 17).Value = 30

# Row 13
$ws.Cells.Item(13, 4).Value = 44757
$ws.Cells.Item(13, 8).Value = "Argentina(o)"
$ws.Cells.Item(13, 9).Value = "Primera"
$ws.Cells.Item(13, 10).Value = 80
$ws.Cells.Item(13, 11).Value = 18000
$ws.Cells.Item(13, 12).Value = 18000
$ws.Cells.Item(13, 13).Value = 18000
$ws.Cells.Item(13, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(13, 16).Value = 360
$ws.Cells.Item(13, 17).Value = 50

# Row 14
$ws.Cells.Item(14, 4).Value = 44757
$ws.Cells.Item(14, 8).Value = "Española"
$ws.Cells.Item(14, 9).Value = "Primera"
$ws.Cells.Item(14, 10).Value = 70
$ws.Cells.Item(14, 11).Value = 22000
$ws.Cells.Item(14, 12).Value = 22000
$ws.Cells.Item(14, 13).Value = 22000
$ws.Cells.Item(14, 14).Value = "`$/caja 30 unidades"
$ws.Cells.Item(14, 16).Value = 733
$ws.Cells.Item(14, 17).Value = 30

# Row 15
$ws.Cells.Item(15, 4).Value = 44488
$ws.Cells.Item(15, 8).Value = "Madrigal"
$ws.Cells.Item(15, 9).Value = "Primera"
$ws.Cells.Item(15, 10).Value = 120
$ws.Cells.Item(15, 11).Value = 12000
$ws.Cells.Item(15, 12).Value = 12000
$ws.Cells.Item(15, 13).Value = 12000
$ws.Cells.Item(15, 14).Value = "`$/caja 40 unidades"
$ws.Cells.Item(15, 16).Value = 300
$ws.Cells.Item(15, 17).Value = 40

# Row 16
$ws.Cells.Item(16, 4).Value = 44418
$ws.Cells.Item(16, 8).Value = "Española"
$ws.Cells.Item(16, 9).Value = "Primera"
$ws.Cells.Item(16, 10).Value = 80
$ws.Cells.Item(16, 11).Value = 16000
$ws.Cells.Item(16, 12).Value = 16000
$ws.Cells.Item(16, 13).Value = 16000
$ws.Cells.Item(16, 14).Value = "`$/caja 30 unidades"
$ws.Cells.Item(16, 16).Value = 533
$ws.Cells.Item(16, 17).Value = 30

# Row 17
$ws.Cells.Item(17, 4).Value = 44775
$ws.Cells.Item(17, 8).Value = "Madrigal"
$ws.Cells.Item(17, 9).Value = "Primera"
$ws.Cells.Item(17, 10).Value = 100
$ws.Cells.Item(17, 11).Value = 16000
$ws.Cells.Item(17, 12).Value = 17000
$ws.Cells.Item(17, 13).Value = 16500
$ws.Cells.Item(17, 14).Value = "`$/caja 40 unidades"
$ws.Cells.Item(17, 16).Value = 412
$ws.Cells.Item(17, 17).Value = 40

# Row 18
$ws.Cells.Item(18, 4).Value = 44386
$ws.Cells.Item(18, 8).Value = "Española"
$ws.Cells.Item(18, 9).Value = "Primera"
$ws.Cells.Item(18, 10).Value = 30
$ws.Cells.Item(18, 11).Value = 15000
$ws.Cells.Item(18, 12).Value = 15000
$ws.Cells.Item(18, 13).Value = 15000
$ws.Cells.Item(18, 14).Value = "`$/caja 30 unidades"
$ws.Cells.Item(18, 16).Value = 500
$ws.Cells.Item(18, 17).Value = 30

# Row 19
$ws.Cells.Item(19, 4).Value = 44778
$ws.Cells.Item(19, 8).Value = "Madrigal"
$ws.Cells.Item(19, 9).Value = "Primera"
$ws.Cells.Item(19, 10).Value = 160
$ws.Cells.Item(19, 11).Value = 15000
$ws.Cells.Item(19, 12).Value = 16000
$ws.Cells.Item(19, 13).Value = 15500
$ws.Cells.Item(19, 14).Value = "`$/caja 40 unidades"
$ws.Cells.Item(19, 16).Value = 388
$ws.Cells.Item(19, 17).Value = 40

# Row 20
$ws.Cells.Item(20, 4).Value = 44421
$ws.Cells.Item(20, 8).Value = "Española"
$ws.Cells.Item(20, 9).Value = "Primera"
$ws.Cells.Item(20, 10).Value = 80
$ws.Cells.Item(20, 11).Value = 16500
$ws.Cells.Item(20, 12).Value = 16500
$ws.Cells.Item(20, 13).Value = 16500
$ws.Cells.Item(20, 14).Value = "`$/caja 30 unidades"
$ws.Cells.Item(20, 16).Value = 550
$ws.Cells.Item(20, 17).Value = 30

# Row 21
$ws.Cells.Item(21, 4).Value = 44771
$ws.Cells.Item(21, 8).Value = "Madrigal"
$ws.Cells.Item(21, 9).Value = "Primera"
$ws.Cells.Item(21, 10).Value = 90
$ws.Cells.Item(21, 11).Value = 16000
$ws.Cells.Item(21, 12).Value = 16000
$ws.Cells.Item(21, 13).Value = 16000
$ws.Cells.Item(21, 14).Value = "`$/caja 40 unidades"
$ws.Cells.Item(21, 16).Value = 400
$ws.Cells.Item(21, 17).Value = 40

# Row 22
$ws.Cells.Item(22, 4).Value = 44782
$ws.Cells.Item(22, 8).Value = "Madrigal"
$ws.Cells.Item(22, 9).Value = "Primera"
$ws.Cells.Item(22, 10).Value = 90
$ws.Cells.Item(22, 11).Value = 15000
$ws.Cells.Item(22, 12).Value = 15000
$ws.Cells.Item(22, 13).Value = 15000
$ws.Cells.Item(22, 14).Value = "`$/caja 40 unidades"
$ws.Cells.Item(22, 16).Value = 375
$ws.Cells.Item(22, 17).Value = 40

# Row 23
$ws.Cells.Item(23, 4).Value = 44484
$ws.Cells.Item(23, 8).Value = "Madrigal"
$ws.Cells.Item(23, 9).Value = "Primera"
$ws.Cells.Item(23, 10).Value = 110
$ws.Cells.Item(23, 11).Value = 11000
$ws.Cells.Item(23, 12).Value = 11000
$ws.Cells.Item(23, 13).Value = 11000
$ws.Cells.Item(23, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(23, 16).Value = 220
$ws.Cells.Item(23, 17).Value = 50
